$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rename header row: "<col>_old" -> "<col>_FV2404", "<col>_new" -> "<col>_FV2410"
# ------------------------------------------------------------------
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value2 = ($baseNames[$i] + "_FV2404")
}

$ws.Cells.Item(1, 11).Value2 = "diff"

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value2 = ($baseNames[$i] + "_FV2410")
}

# ------------------------------------------------------------------
# 2) Turn the used range into an Excel Table ("Table1") without
#    disturbing the existing header formatting (bold/fill/border/center/wrap).
#    ListObjects.Add() bakes the pre-existing header format into a dxf
#    (headerRowDxfId) the first time it sees an already-formatted header
#    row, so we stash a copy of that formatting on a scratch cell, clear
#    the header, create the table, and paste the formatting back in one
#    shot so the original style index is reused untouched.
# ------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("W1")

$ws.Range("A1").Copy($scratch) | Out-Null
$headerRange.ClearFormats() | Out-Null

$dataRange = $ws.Range("A1:U91")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

$scratch.Copy() | Out-Null
$headerRange.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$scratch.Clear() | Out-Null
$excel.CutCopyMode = $false

$tbl.TableStyle = ""

# ------------------------------------------------------------------
# 3) Freeze the header row (split above row 2).
# ------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Host "edit.ps1 completed"
